# Add new columns I ("I0") and J ("IF") to Sheet1, mirroring existing
# header style from column H, and populate data rows 2-66.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): I1 = "I0", J1 = "IF" ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting from the existing header cell H1 (bold, border,
# centered/top-aligned) onto the two new header cells so they match the
# rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# --- Data rows (rows 2-66): numeric values for I0 and IF ---
$I0Values = @(7,6,8,8,7,7,5,6,6,7,6,7,5,10,6,1,5,8,7,7,7,6,8,7,7,6,7,6,7,5,7,6,6,5,7,8,8,7,7,7,6,6,6,7,7,6,7,8,9,8,8,7,6,7,5,6,7,6,5,8,6,1,1,5,1)
$IFValues = @(7,7,9,9,7,7,6,6,6,7,6,8,6,10,6,1,5,8,7,7,7,7,8,8,7,6,7,7,7,7,8,8,6,6,7,8,9,8,8,7,6,7,8,9,8,6,7,8,9,9,8,8,8,8,7,7,7,8,6,8,7,4,5,7,2)

for ($i = 0; $i -lt $I0Values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 9).Value  = $I0Values[$i]
    $ws.Cells.Item($row, 10).Value = $IFValues[$i]
}

Write-Output "I0 and IF columns added"
